# Replace the single M2Doc field (fldChar begin/end + instrText runs) with
# plain-text runs wrapped in "{" / "}" delimiters, as produced by the new
# TokenIteratorFieldRewriterSplit parser.
#
# Before: <w:r><w:fldChar begin/></w:r>
#         <w:r><w:instrText> </w:instrText></w:r>        (leading space)
#         ... several <w:instrText> runs (code text, with a bookmark in the
#             middle) ...
#         <w:r><w:instrText> </w:instrText></w:r>        (trailing space)
#         <w:r><w:fldChar end/></w:r>
#
# After:  <w:r><w:t>{</w:t></w:r>
#         ... the same runs, but <w:instrText> -> <w:t> ...
#         <w:r><w:t>}</w:t></w:r>
# (the leading/trailing single-space runs become "{" and "}" respectively,
# and the begin/end fldChar runs are dropped entirely)

$d = $word.ActiveDocument

# Locate the paragraph that owns the (only) field in the document.
$field = $d.Fields.Item(1)
$fieldStart = $field.Code.Start

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if (($fieldStart -ge $p.Range.Start) -and ($fieldStart -lt $p.Range.End)) {
        $target = $p
        break
    }
}

$full = $target.Range
$contentRange = $d.Range($full.Start, $full.End - 1)

$newRunsXml = '<w:r><w:t>{</w:t></w:r>' + `
    '<w:r><w:t>m</w:t></w:r>' + `
    '<w:r><w:t>:</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' + `
    '<w:r><w:t>self.na</w:t></w:r>' + `
    '<w:r><w:t>me + &apos;</w:t></w:r>' + `
    '<w:r><w:t>\n</w:t></w:r>' + `
    '<w:r><w:t>\n![](../images/logo_M2Doc.png</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:t>)</w:t></w:r>' + `
    '<w:r><w:t>&apos;).from</w:t></w:r>' + `
    '<w:r><w:t>Markdown</w:t></w:r>' + `
    '<w:r><w:t>String(</w:t></w:r>' + `
    '<w:r><w:t>&apos;https://www.m2doc.org/tests/&apos;</w:t></w:r>' + `
    '<w:r><w:t>)</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p>' + $newRunsXml + '</w:p></w:body>' + `
    '</w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$contentRange.InsertXML($packageXml) | Out-Null
